# Word COM-interop script implementing the "Use cases.docx" edit:
#   1. Hashem row / "Program used" cell: split "...Word, Github" so that
#      "Github" is wrapped in spell-check proofErr markers.
#   2. Ismail row / "Program used" cell: append ", Docker" after "Github".
#   3. Ismail row / "Tasks done" cell: rewrite the paragraph to add the new
#      sentences about the Dockerfile, dropping the old grammar proofErr
#      markers around "helped out".
#
# We use Range.InsertXML with a full WordProcessingML "flat OPC" package
# wrapper (the same shape Range.WordOpenXML emits) because it lets us
# control run-level markup precisely, including <w:proofErr/> elements that
# have no Range.Text-level representation.

$d = $word.ActiveDocument

function New-FlatOpcPackage([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Set-CellParagraphXml($cell, [string]$paragraphXml) {
    # cell.Range spans the cell's text *plus* the trailing end-of-cell mark;
    # trimming the last character keeps that mark untouched so the table
    # structure stays intact, and lets InsertXML fully replace the
    # paragraph(s) that make up the cell's visible content.
    $full = $cell.Range
    $target = $d.Range($full.Start, $full.End - 1)
    $target.InsertXML((New-FlatOpcPackage $paragraphXml))
}

$table = $d.Tables.Item(1)

# --- 1. Hashem / Program used ------------------------------------------
$hashemCell = $table.Rows.Item(3).Cells.Item(2)
$target1 = '<w:p w14:paraId="3F1A2FD4" w14:textId="70F470C9" w:rsidR="000B37D3" w:rsidRDefault="000B37D3">' +
    '<w:r><w:t xml:space="preserve">Visual Studio Code, Docker, Node.js, Word, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Github</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Set-CellParagraphXml $hashemCell $target1

# --- 2. Ismail / Program used -------------------------------------------
$ismailProgCell = $table.Rows.Item(5).Cells.Item(2)
$target2 = '<w:p w14:paraId="338C8786" w14:textId="449C8AAD" w:rsidR="000B37D3" w:rsidRDefault="000B37D3">' +
    '<w:r><w:t xml:space="preserve">Zube, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Github</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, Docker</w:t></w:r>' +
    '</w:p>'
Set-CellParagraphXml $ismailProgCell $target2

# --- 3. Ismail / Tasks done ----------------------------------------------
$ismailTaskCell = $table.Rows.Item(5).Cells.Item(3)
$bodyRpr = '<w:rPr><w:rFonts w:ascii="UICTFontTextStyleBody" w:hAnsi="UICTFontTextStyleBody"/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
$trailRpr = '<w:rPr><w:rFonts w:ascii="UICTFontTextStyleBody" w:hAnsi="UICTFontTextStyleBody"/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
$target3 = '<w:p w14:paraId="3137EBF5" w14:textId="07BD8969" w:rsidR="000B37D3" w:rsidRPr="000B37D3" w:rsidRDefault="000B37D3">' +
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r>' + $bodyRpr + '<w:t xml:space="preserve">He integrated the project with Zube.io. </w:t></w:r>' +
    '<w:r>' + $bodyRpr + '<w:t xml:space="preserve">Ismail also created the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $bodyRpr + '<w:t>Dockerfile</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $bodyRpr + '<w:t xml:space="preserve">. Wrote the code in the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $bodyRpr + '<w:t>dockerfile</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $bodyRpr + '<w:t xml:space="preserve">. </w:t></w:r>' +
    '<w:r>' + $bodyRpr + '<w:t>He also helped out with incorporating the Kanban/Project board.</w:t></w:r>' +
    '<w:r>' + $trailRpr + '<w:t> </w:t></w:r>' +
    '</w:p>'
Set-CellParagraphXml $ismailTaskCell $target3

Write-Output "done"
